# edit.ps1
# Applies the NextGenFwys ModelRuns.xlsx update:
#  - Renames category "NGF" -> "NextGenFwys" for all NGF rows (col A)
#  - Renames a handful of run names (col C) to reflect the new naming scheme
#  - Inserts three additional sensitivity/toll-level rows before the former
#    last (Blueprint) row, and appends one more "BlueprintSegmentedTest" row
#    after it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all_runs")

# --- 1. Column A: "NGF" -> "NextGenFwys" for the project rows that use it ---
foreach ($r in 3..13) {
    if ($ws.Cells.Item($r, 1).Value -eq "NGF") {
        $ws.Cells.Item($r, 1).Value = "NextGenFwys"
    }
}

# --- 2. Rename a few existing run names (column C) ---
$ws.Cells.Item(8, 3).Value = "2035_TM152_NGF_NP01"
$ws.Cells.Item(9, 3).Value = "2035_TM152_NGF_NP02_Blueprint_00_SensDiscount_01"
$ws.Cells.Item(10, 3).Value = "2035_TM152_NGF_NP02_Blueprint_00_SensDiscount_02"
$ws.Cells.Item(11, 3).Value = "2035_TM152_NGF_NP02_Blueprint_00_SensDiscount_03"
$ws.Cells.Item(12, 3).Value = "2035_TM152_NGF_NP02_Blueprint_00_SensDiscount_04"
$ws.Cells.Item(13, 3).Value = "2035_TM152_NGF_NP02_Blueprint_00_SensDiscount_05"

# --- 3. Insert 3 new rows before the old row 14 (Blueprint row), copying
#        row 13's formatting, then fill in the new data ---
$ws.Rows.Item(14).Resize(3,1).EntireRow.Insert()

$ws.Rows.Item(13).EntireRow.Copy()
$ws.Rows.Item(14).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Rows.Item(15).PasteSpecial(-4122) | Out-Null
$ws.Rows.Item(16).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$rowsData = @(
    @{ Row = 14; C = "2035_TM152_NGF_NP02_Blueprint_00_SensExtent_01" },
    @{ Row = 15; C = "2035_TM152_NGF_NP02_Blueprint_00_TollLevel_01" },
    @{ Row = 16; C = "2035_TM152_NGF_NP02_Blueprint_00_TollLevel_02" }
)

foreach ($rd in $rowsData) {
    $r = $rd.Row
    $ws.Cells.Item($r, 1).Value = "NextGenFwys"
    $ws.Cells.Item($r, 2).Value = 2035
    $ws.Cells.Item($r, 3).Value = $rd.C
    $ws.Cells.Item($r, 4).Value = "NGF"
    $ws.Cells.Item($r, 5).Value = "Sensitivity Test"
    $ws.Cells.Item($r, 6).Value = '"Final Blueprint runs\Final Blueprint (s24)\BAUS v2.25 - FINAL VERSION"'
    $ws.Cells.Item($r, 7).Value = "run182"
    $ws.Cells.Item($r, 8).Value = "current"
}

# --- 4. Append a new row (18) after the Blueprint row (now row 17), copying
#        row 17's formatting (the hyperlinked row) ---
$ws.Rows.Item(17).EntireRow.Copy()
$ws.Rows.Item(18).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(18, 1).Value = "NextGenFwys"
$ws.Cells.Item(18, 2).Value = 2035
$ws.Cells.Item(18, 3).Value = "2035_TM152_NGF_NP02_Blueprint_00_BlueprintSegmentedTest"
$ws.Cells.Item(18, 4).Value = "NGF"
$ws.Cells.Item(18, 5).Value = "Blueprint"
$ws.Cells.Item(18, 6).Value = '"Final Blueprint runs\Final Blueprint (s24)\BAUS v2.25 - FINAL VERSION"'
$ws.Cells.Item(18, 7).Value = "run182"
$ws.Cells.Item(18, 8).Value = "current"
$ws.Cells.Item(18, 9).Value = ""
$ws.Cells.Item(18, 10).Value = ""

# --- 5. Column formatting / view tweaks ---
$ws.Columns.Item(3).ColumnWidth = 50.3984375

$ws.Range("C25").Select()
